$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1118
$ws1.Range("F3").Value = 633
$ws1.Range("F4").Value = 0
$ws1.Range("F6").Value = 513
$ws1.Range("F7").Value = 8988
$ws1.Range("F8").Value = 0
$ws1.Range("F11").Value = 598

# Sheet "演出" (sheet2): update F column
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 0

# Sheet "全部类型" (sheet4): update F column
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1118
$ws4.Range("F3").Value = 633
$ws4.Range("F5").Value = 0
$ws4.Range("F7").Value = 4905
$ws4.Range("F8").Value = 513
$ws4.Range("F9").Value = 7
$ws4.Range("F10").Value = 8988
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F16").Value = 598
$ws4.Range("F17").Value = 65
